$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.1682750951017979
$ws.Range("B2").Value = -0.1682750951017979
$ws.Range("C2").Value = -0.1682750951017979
$ws.Range("D2").Value = -0.1682750951017979
$ws.Range("E2").Value = -0.1682750951017979
$ws.Range("F2").Value = -0.1682750951017979
$ws.Range("G2").Value = -0.1682750951017979
$ws.Range("H2").Value = -0.1678866623857618
$ws.Range("I2").Value = -0.1748153949307012
$ws.Range("J2").Value = -0.170453939667295
$ws.Range("K2").Value = -0.1727477127267293
$ws.Range("L2").Value = -0.1682750951017979
$ws.Range("M2").Value = -0.178677714881075
$ws.Range("N2").Value = -0.1682750951017979
$ws.Range("O2").Value = -0.174371191312673
$ws.Range("P2").Value = -0.1682750951017979
$ws.Range("Q2").Value = -0.1622525240592109
$ws.Range("R2").Value = -0.1682750951017979
$ws.Range("S2").Value = -0.1622525240592109
$ws.Range("T2").Value = -0.1682750951017979
